## Fixed update to excel issue
## - Rename "Requested quantity" headers to source-specific names
## - Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

# Rename headers
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add new "PO Forecast" sheet after the last existing sheet
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Match header styling (bold, centered, bordered) from an existing sheet's header
$wsWeekly.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$ws3.Range("A2").Value = 45137.99999999999
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = -265.8269485012258
$ws3.Range("D2").Value = 65.54797327841125
$ws3.Range("A3").Value = 45144.99999999999
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = -241.7686739057561
$ws3.Range("D3").Value = 88.02394547933854
$ws3.Range("A4").Value = 45151.99999999999
$ws3.Range("B4").Value = 0
$ws3.Range("C4").Value = -231.1193107896818
$ws3.Range("D4").Value = 105.0514776843997
$ws3.Range("A5").Value = 45158.99999999999
$ws3.Range("B5").Value = 0
$ws3.Range("C5").Value = -197.0066266754478
$ws3.Range("D5").Value = 120.3500004420433
$ws3.Range("A6").Value = 45186.99999999999
$ws3.Range("B6").Value = 50
$ws3.Range("C6").Value = -104.7600718371033
$ws3.Range("D6").Value = 197.1847716712437
$ws3.Range("A7").Value = 45207.99999999999
$ws3.Range("B7").Value = 113
$ws3.Range("C7").Value = -46.11588207558689
$ws3.Range("D7").Value = 277.8414085411758
$ws3.Range("A8").Value = 45214.99999999999
$ws3.Range("B8").Value = 134
$ws3.Range("C8").Value = -31.51456434206747
$ws3.Range("D8").Value = 297.1135971015117
$ws3.Range("A9").Value = 45221.99999999999
$ws3.Range("B9").Value = 155
$ws3.Range("C9").Value = -5.98860406231353
$ws3.Range("D9").Value = 312.9476249603375
$ws3.Range("A10").Value = 45228.99999999999
$ws3.Range("B10").Value = 176
$ws3.Range("C10").Value = 7.405744984650082
$ws3.Range("D10").Value = 346.2283558114841
$ws3.Range("A11").Value = 45235.99999999999
$ws3.Range("B11").Value = 197
$ws3.Range("C11").Value = 30.44869702523624
$ws3.Range("D11").Value = 352.1118849297521
$ws3.Range("A12").Value = 45242.99999999999
$ws3.Range("B12").Value = 218
$ws3.Range("C12").Value = 56.26829269122157
$ws3.Range("D12").Value = 379.3730459191249
$ws3.Range("A13").Value = 45249.99999999999
$ws3.Range("B13").Value = 239
$ws3.Range("C13").Value = 79.15942111906421
$ws3.Range("D13").Value = 403.8232705902613
$ws3.Range("A14").Value = 45263.99999999999
$ws3.Range("B14").Value = 281
$ws3.Range("C14").Value = 128.1381026595836
$ws3.Range("D14").Value = 439.6441002460469
$ws3.Range("A15").Value = 45298.99999999999
$ws3.Range("B15").Value = 386
$ws3.Range("C15").Value = 220.9337730557192
$ws3.Range("D15").Value = 544.5453460344065
$ws3.Range("A16").Value = 45326.99999999999
$ws3.Range("B16").Value = 470
$ws3.Range("C16").Value = 299.5489285020234
$ws3.Range("D16").Value = 644.4560647789617
$ws3.Range("A17").Value = 45333.99999999999
$ws3.Range("B17").Value = 491
$ws3.Range("C17").Value = 324.9774017950258
$ws3.Range("D17").Value = 661.6027700264656
$ws3.Range("A18").Value = 45340.99999999999
$ws3.Range("B18").Value = 513
$ws3.Range("C18").Value = 350.4036466941164
$ws3.Range("D18").Value = 672.9597983192131
$ws3.Range("A19").Value = 45347.99999999999
$ws3.Range("B19").Value = 534
$ws3.Range("C19").Value = 370.0859107645706
$ws3.Range("D19").Value = 694.4246326925141
$ws3.Range("A20").Value = 45354.99999999999
$ws3.Range("B20").Value = 555
$ws3.Range("C20").Value = 385.6993167358921
$ws3.Range("D20").Value = 712.1861776772017
$ws3.Range("A21").Value = 45361.99999999999
$ws3.Range("B21").Value = 576
$ws3.Range("C21").Value = 408.8856595295068
$ws3.Range("D21").Value = 735.7948680386878
$ws3.Range("A22").Value = 45368.99999999999
$ws3.Range("B22").Value = 597
$ws3.Range("C22").Value = 419.5575047739869
$ws3.Range("D22").Value = 764.6910982894377
$ws3.Range("A23").Value = 45375.99999999999
$ws3.Range("B23").Value = 618
$ws3.Range("C23").Value = 468.5977580910555
$ws3.Range("D23").Value = 783.2631069534166
$ws3.Range("A24").Value = 45382.99999999999
$ws3.Range("B24").Value = 639
$ws3.Range("C24").Value = 479.2923069465805
$ws3.Range("D24").Value = 813.2195925389061
$ws3.Range("A25").Value = 45389.99999999999
$ws3.Range("B25").Value = 660
$ws3.Range("C25").Value = 497.79579361214
$ws3.Range("D25").Value = 834.2404067927323

# Match date styling for the date column from an existing sheet's date cell
$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A25").PasteSpecial(-4122)

# Restore original selection on the first sheet so active sheet/selection match source
$wsWeekly.Select() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
